$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IT")
$ws.Range("B2").Value = 2019
